# Add "Diode 1N4007" to the purchase list (row 30), pushing existing
# rows 30-55 down to 31-56, and refresh the autofilter / filter-database
# defined name + selection to match.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new blank row at row 30 (shifts old rows 30-55 -> 31-56).
$null = $ws.Rows("30:30").Insert()

# Populate the new row with the diode purchase info.
$ws.Range("A30").Value = "Diode 1N4007"
$ws.Range("B30").Value = "Amazon"
$ws.Range("C30").Value = 125
$ws.Range("D30").Value = 1
$ws.Range("E30").Value = 5.99
$ws.Range("F30").Value = 1
$ws.Range("G30").Value = "https://www.amazon.com/BOJACK-Rectifier-IN4007-Electronic-Silicon/dp/B07Q6J9TNW"
$ws.Range("H30").Value = "1/ solution"

# Extend the autofilter over the new data range (A1:H56).
$ws.AutoFilterMode = $false
$null = $ws.Range("A1:H56").AutoFilter()

# Update the _FilterDatabase defined name to match the new range.
foreach ($n in $wb.Names) {
  if ($n.Name -like "*_FilterDatabase*") {
    $n.RefersTo = "=Sheet1!`$A`$1:`$H`$56"
  }
}

# Match the author's final selection/scroll state.
$null = $ws.Range("A30:G30").Select()
